$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row at position 13 (old rows 13-23 shift down to 14-24).
$ws.Range("A13").EntireRow.Insert()

# The inserted row inherits column A's formatting from the row above; clear it so
# row 13 ends up with no A cell at all (matching the target layout).
$ws.Range("A13").Clear()

# Give the new B13:C13 cells the same "plain" (non-header) formatting used by the
# other two-column data rows (e.g. B10:C10), then fill them in.
$ws.Range("B10:C10").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$ws.Range("B13").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("C13").Value = "984972 - Hugo Ricardo Zschommler Sandim"

# 2. Fix up the various answer cells with their correct text.

# Objetivos:
$ws.Range("B10").Value = "Apresentar os principais conceitos sobre as transformações de fases em materiais metálicos, poliméricos e cerâmicos abrangendo transformações difusionais e não-difusionais, a conceituação sobre nucleação e crescimento (aspectos energéticos) e sua relação com problemas práticos encontrados nas indústrias de processamento e de transformação de materiais."
$ws.Range("C10").Value = "Apresentar os principais conceitos sobre as transformações de fases em materiais metálicos, poliméricos e cerâmicos abrangendo transformações difusionais e não-difusionais, a conceituação sobre nucleação e crescimento (aspectos energéticos) e sua relação com problemas práticos encontrados nas indústrias de processamento e de transformação de materiais."

# Programa resumido: (now row 14)
$ws.Range("B14").Value = "Difusão no estado sólido.Difusão em materiais não-metálicos (sólidos iônicos e polímeros).Recuperação, recristalização e crescimento de grão.Solidificação.Precipitação.Cinética de transformação no sistema Fe-C e em ligas não-ferrosas.Transformação de fases em vidros e cerâmicas.Transformação de fases em polímeros.Atividade experimental."
$ws.Range("C14").Value = "Difusão no estado sólido.Difusão em materiais não-metálicos (sólidos iônicos e polímeros).Recuperação, recristalização e crescimento de grão.Solidificação.Precipitação.Cinética de transformação no sistema Fe-C e em ligas não-ferrosas.Transformação de fases em vidros e cerâmicas.Transformação de fases em polímeros.Atividade experimental."

# Programa: (now row 16)
$ws.Range("B16").Value = "Introdução à difusão no estado sólido. Coeficiente de difusão. Leis de Fick. Difusão em soluções diluídas e na presença de um gradiente de concentração. Efeito Kirkendall.- Apresentar os fundamentos teóricos pertinentes à transformação de fases em materiais metálicos, cerâmicos e poliméricos.- Apresentar os conceitos fundamentais associados à nucleação, ao crescimento e à cinética de transformação de fases.- Descrições detalhadas de microestruturas fundidas e tratadas termicamente. Aspectos morfológicos relevantes.- Descrição das principais transformações de fase no estado sólido no sistema Fe-C e em algumas ligas não-ferrosas. Curvas TTT e CCT (TRC).- Estudar a transformação de fases durante o processamento termomecânico de metais e ligas.- Realização de prática experimental versando sobre tópicos da ementa."
$ws.Range("C16").Value = "Introdução à difusão no estado sólido. Coeficiente de difusão. Leis de Fick. Difusão em soluções diluídas e na presença de um gradiente de concentração. Efeito Kirkendall.- Apresentar os fundamentos teóricos pertinentes à transformação de fases em materiais metálicos, cerâmicos e poliméricos.- Apresentar os conceitos fundamentais associados à nucleação, ao crescimento e à cinética de transformação de fases.- Descrições detalhadas de microestruturas fundidas e tratadas termicamente. Aspectos morfológicos relevantes.- Descrição das principais transformações de fase no estado sólido no sistema Fe-C e em algumas ligas não-ferrosas. Curvas TTT e CCT (TRC).- Estudar a transformação de fases durante o processamento termomecânico de metais e ligas.- Realização de prática experimental versando sobre tópicos da ementa."

# Método: (now row 19)
$ws.Range("B19").Value = "O aluno será avaliado ao longo do semestre por duas avaliações escritas (P1 e P2) e com pesos iguais."
$ws.Range("C19").Value = "O aluno será avaliado ao longo do semestre por duas avaliações escritas (P1 e P2) e com pesos iguais."

# Critério: (now row 20)
$ws.Range("B20").Value = "Nota Final NF = [P1 + P2]/2"
$ws.Range("C20").Value = "Nota Final NF = [P1 + P2]/2"

# Norma de recuperação: (now row 21)
$ws.Range("B21").Value = "Para a recuperação será realizada uma prova escrita (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2"
$ws.Range("C21").Value = "Para a recuperação será realizada uma prova escrita (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2"

# Bibliografia: (now row 22)
$ws.Range("B22").Value = "1 - Diffusion in solids. P.G. Shewmon, McGraw-Hill, 1963.2 - Phase transformation in metals. P.G. Shewmon, McGraw-Hill, 1969.3 - Recrystallization and related annealing phenomena. F.J. Humphreys and M. Hatherly, Pergamon, 1996.4 - Principles of solidification. B. Chalmers, Robert E. Krieger, 2nd. ed., 1977.5 - Precipitation hardening. A. Kelly, Pergamon, 1963.6 - Particle strengthening of metals and alloys. E. Nembach, John Wiley & Sons, 1997.7 - Propriedades dos materiais cerâmicos. L.H. Van Vlack, Edgard Blücher, 1973.8 - Textbook of polymer science. F.W. Billmeyer Jr., John Wiley & Sons, 1962.9 - Worked examples in the kinetics and thermodynamics of phase transformations. E.A. Wilson, The Institution of Metallurgists, s.d."
$ws.Range("C22").Value = "1 - Diffusion in solids. P.G. Shewmon, McGraw-Hill, 1963.2 - Phase transformation in metals. P.G. Shewmon, McGraw-Hill, 1969.3 - Recrystallization and related annealing phenomena. F.J. Humphreys and M. Hatherly, Pergamon, 1996.4 - Principles of solidification. B. Chalmers, Robert E. Krieger, 2nd. ed., 1977.5 - Precipitation hardening. A. Kelly, Pergamon, 1963.6 - Particle strengthening of metals and alloys. E. Nembach, John Wiley & Sons, 1997.7 - Propriedades dos materiais cerâmicos. L.H. Van Vlack, Edgard Blücher, 1973.8 - Textbook of polymer science. F.W. Billmeyer Jr., John Wiley & Sons, 1962.9 - Worked examples in the kinetics and thermodynamics of phase transformations. E.A. Wilson, The Institution of Metallurgists, s.d."

# 3. Column widths: narrow the "col min=1 max=2" definition down to just column 1,
# since column 2 already has its own explicit (identical-width) definition below it.
$ws.Columns.Item(1).ColumnWidth = 30.7109375
